$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the previously-missing data for row 16 (Problem 3, h_ignore_preconditions)
$ws.Range("C16").Value = 408
$ws.Range("D16").Value = 410
$ws.Range("E16").Value = 3758
$ws.Range("F16").Value = 12
$ws.Range("G16").Value = 586.231951499042

# Update the selected/active cell shown in the saved view
$ws.Range("D17").Select()
